$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the ROC (column B) value for every ticker currently on the sheet.
#    Values are assigned by the row's CURRENT ticker identity (i.e. the data
#    update happens before the table gets re-sorted below).
$ws.Range("B2").Value = 37.02    # IBIT
$ws.Range("B3").Value = 22.81    # GLD
$ws.Range("B4").Value = 7.06     # FXI
$ws.Range("B5").Value = 2.77     # XLP
$ws.Range("B6").Value = 1.94     # XLF
$ws.Range("B7").Value = -1.46    # XLU
$ws.Range("B8").Value = -3.69    # SLV
$ws.Range("B9").Value = -2.95    # UUP
$ws.Range("B10").Value = -4.15   # TLT
$ws.Range("B11").Value = -4.52   # MTUM
$ws.Range("B12").Value = -6.55   # USO
$ws.Range("B13").Value = -6.86   # XLY
$ws.Range("B14").Value = -6.7    # XLRE
$ws.Range("B15").Value = -7.21   # XLI
$ws.Range("B16").Value = -7.17   # GMF
$ws.Range("B17").Value = -7.86   # XLE
$ws.Range("B18").Value = -7.54   # RSP
$ws.Range("B19").Value = -8.36   # SPY
$ws.Range("B20").Value = -8.44   # SPYV
$ws.Range("B21").Value = -8.34   # XLV
$ws.Range("B22").Value = -9.21   # SPYG
$ws.Range("B23").Value = -9.99   # QQQ
$ws.Range("B24").Value = -12.52  # XLB
$ws.Range("B25").Value = -13.18  # MOAT
$ws.Range("B26").Value = -13.51  # IWN
$ws.Range("B27").Value = -14.37  # IWO
$ws.Range("B28").Value = -16.14  # XLK

# 2) Re-rank the table: sort rows 2-28 by the (now updated) ROC column,
#    descending, carrying each ticker's full row (and its per-row
#    formatting) along with it.
$rng = $ws.Range("A2:F28")
$key = $ws.Range("B2:B28")
$rng.Sort($key, 2)

# 3) One row's "Mensual" indicator (column C) was independently re-flagged
#    by the source data refresh and doesn't follow purely from the
#    re-sort: after sorting, XLRE (new row 13) keeps its existing
#    highlight, while XLY (new row 14) picks up the "Mensual" flag that
#    C16 (style s=8) already shows. Copy that formatting over.
$ws.Range("C16").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null
